# #165 point score calculation logic improvement
# Add a new match (Abu Dhabi & Dubai fixtures) as a new first sheet ("Sheet3"),
# pushing the existing Sheet1/Sheet2 tabs after it.

$wb = $excel.ActiveWorkbook

# Add a brand-new worksheet; Excel inserts it before the active sheet and
# assigns it the next free sheet name ("Sheet3"), matching the target
# workbook's sheet order of Sheet3, Sheet1, Sheet2.
$newSheet = $wb.Worksheets.Add()

# Match the new sheet's column widths to the authored layout.
$newSheet.Columns.Item(1).ColumnWidth = 26.6640625
$newSheet.Columns.Item(2).ColumnWidth = 24.83203125
$newSheet.Columns.Item(3).ColumnWidth = 21.83203125
$newSheet.Columns.Item(5).ColumnWidth = 31.33203125

# Header row
$newSheet.Range("A1").Value = "Home_Team"
$newSheet.Range("B1").Value = "Away_Team"
$newSheet.Range("C1").Value = "Tournament"
$newSheet.Range("D1").Value = "Venue"
$newSheet.Range("E1").Value = "matchTime"

# First new fixture: Mumbai Indians vs Chennai Super Kings, Abu Dhabi
$newSheet.Range("A2").Value = "Mumbai Indians"
$newSheet.Range("B2").Value = "Chennai Super Kings"
$newSheet.Range("C2").Value = "IPL-20"
$newSheet.Range("D2").Value = "ABU DHABI"
$newSheet.Range("E2").Value = "2020-09-06T10:00:00.000Z"

# Second new fixture: Delhi Capitals vs Kings XI Punjab, Dubai
$newSheet.Range("C3").Value = "IPL-20"
$newSheet.Range("D3").Value = "DUBAI"
$newSheet.Range("E3").Value = "2020-09-06T13:30:00.000Z"
$newSheet.Range("B3").Value = "Kings XI Punjab"
$newSheet.Range("A3").Value = "Delhi Capitals"

# Restore the previous selection on the (now second) original schedule sheet.
$sheet2 = $wb.Worksheets.Item("Sheet2")
$sheet2.Range("A4").Select()

# Leave the new sheet active/selected on B2, as authored.
$newSheet.Range("B2").Select()
